$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string label swaps (reorder countries) ---
$ws.Range("A58").Value = "Argelia"
$ws.Range("A59").Value = "Suiza"
$ws.Range("A114").Value = "Namibia"
$ws.Range("A115").Value = "Congo"
$ws.Range("A156").Value = "Aruba"
$ws.Range("A157").Value = "Principado de Andorra"

# --- Updated COVID numbers per row ---
$ws.Range("B4").Value = 5501544
$ws.Range("C4").Value = 25278
$ws.Range("D4").Value = 2879670
$ws.Range("E4").Value = 2449755
$ws.Range("G4").Value = 584
$ws.Range("H4").Value = 172119

$ws.Range("B6").Value = 2588253
$ws.Range("C6").Value = 63031
$ws.Range("D6").Value = 1858983
$ws.Range("E6").Value = 679186
$ws.Range("G6").Value = 950
$ws.Range("H6").Value = 50084

$ws.Range("B21").Value = 248117
$ws.Range("C21").Value = 1256
$ws.Range("D21").Value = 229972
$ws.Range("E21").Value = 12190
$ws.Range("G21").Value = 21
$ws.Range("H21").Value = 5955

$ws.Range("B22").Value = 224360
$ws.Range("C22").Value = 586
$ws.Range("E22").Value = 12520
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 9290

$ws.Range("B33").Value = 92198
$ws.Range("C33").Value = 1118
$ws.Range("D33").Value = 67950
$ws.Range("E33").Value = 23574
$ws.Range("G33").Value = 9
$ws.Range("H33").Value = 674

$ws.Range("B58").Value = 38133
$ws.Range("C58").Value = 469
$ws.Range("D58").Value = 26644
$ws.Range("E58").Value = 10129
$ws.Range("G58").Value = 9
$ws.Range("H58").Value = 1360

$ws.Range("B59").Value = 37924
$ws.Range("C59").Value = 253
$ws.Range("D59").Value = 33200
$ws.Range("E59").Value = 2733
$ws.Range("H59").Value = 1991

$ws.Range("B61").Value = 34528
$ws.Range("C61").Value = 707
$ws.Range("D61").Value = 29328
$ws.Range("E61").Value = 4975
$ws.Range("G61").Value = 5
$ws.Range("H61").Value = 225

$ws.Range("B104").Value = 5679
$ws.Range("C104").Value = 107
$ws.Range("D104").Value = 3208
$ws.Range("E104").Value = 2449

$ws.Range("B114").Value = 3907
$ws.Range("C114").Value = 181
$ws.Range("D114").Value = 2352
$ws.Range("E114").Value = 1520
$ws.Range("G114").Value = 4
$ws.Range("H114").Value = 35

$ws.Range("B115").Value = 3745
$ws.Range("D115").Value = 1625
$ws.Range("E115").Value = 2060
$ws.Range("H115").Value = 60

$ws.Range("B156").Value = 1048
$ws.Range("C156").Value = 75
$ws.Range("D156").Value = 193
$ws.Range("E156").Value = 851
$ws.Range("H156").Value = 4

$ws.Range("B157").Value = 989
$ws.Range("D157").Value = 863
$ws.Range("E157").Value = 73
$ws.Range("H157").Value = 53

# --- Timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Agosto de 2020 a las 20:54"
